# Insert a new weekly record right before the current row 413.
# Excel's native row insert shifts rows 413:449 down to 414:450
# (carrying their values/styles along), which matches the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(413).Insert()

# Populate the newly inserted row 413 with the new weekly record.
$ws.Range("A413").Value = 8
$ws.Range("B413").Value = "Terminal La Palmera de La Serena"
$ws.Range("C413").Value = "Coquimbo"
$ws.Range("D413").Value = 45013
$ws.Range("E413").Value = 4
$ws.Range("F413").Value = 100112003
$ws.Range("G413").Value = "Ajo"
$ws.Range("H413").Value = "Chino"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 440
$ws.Range("K413").Value = 17000
$ws.Range("L413").Value = 18000
$ws.Range("M413").Value = 17500
$ws.Range("N413").Value = "$/caja 10 kilos"
$ws.Range("O413").Value = "China"
$ws.Range("P413").Value = 1750
$ws.Range("Q413").Value = 10
$ws.Range("R413").Value = "Hortaliza"
